# Germany Verbandsliga.xlsx update
# - Swap the match-result data (everything except id/Div/Div Original Name/Date/HomeTeam)
#   between a few row pairs that were re-sequenced at the source (rows 5/7, 46/47, 67/68).
# - Append 5 new fixtures (rows 112-116).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Swap AwayTeam + result/odds columns between paired rows.
#    Columns kept in place: A (id), C (Div), D (Div Original Name), E (Date), F (HomeTeam)
#    Columns swapped: B (match id) and G..AC (AwayTeam through PL_AhUnder)
# ---------------------------------------------------------------------------
function Swap-Rows($rowA, $rowB) {
    $cols = @(2,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29)
    foreach ($col in $cols) {
        $valA = $ws.Cells.Item($rowA, $col).Value2
        $valB = $ws.Cells.Item($rowB, $col).Value2
        $ws.Cells.Item($rowA, $col).Value2 = $valB
        $ws.Cells.Item($rowB, $col).Value2 = $valA
    }
}

Swap-Rows 5 7
Swap-Rows 46 47
Swap-Rows 67 68

# ---------------------------------------------------------------------------
# 2) Append 5 new fixture rows (112-116), copying the row-111 formatting
#    (bold/bordered style on column A, date number format on column E).
# ---------------------------------------------------------------------------
$ws.Range("A111:AC111").Copy()
$ws.Range("A112:AC116").PasteSpecial(-4122)
$excel.CutCopyMode = 0

function Set-Row($r, $values) {
    $arr = New-Object 'object[,]' 1,29
    for ($i = 0; $i -lt 29; $i++) {
        $arr[0,$i] = $values[$i]
    }
    $ws.Range($ws.Cells.Item($r,1), $ws.Cells.Item($r,29)).Value2 = $arr
}

Set-Row 112 @(110,7981940,"Germany Verbandsliga","Germany Verbandsliga",45371.64583333334,"VfB Wissen","Niederrossbach",6,0,"H",1.333,5,6,1.333,5.25,5.75,-1.75,1.975,1.825,4,1.9,1.9,0.333,-1,-1,0.9750000000000001,-1,0.8999999999999999,-1)

Set-Row 113 @(111,7990779,"Germany Verbandsliga","Germany Verbandsliga",45373.625,"BSV HalleAmmendorf","SG RotWeiss Thalheim",1,2,"A",1.909,4,2.9,1.909,4,2.9,-0.5,1.975,1.825,3.25,2.025,1.775,-1,-1,1.9,-1,0.825,-0.5,0.3875)

Set-Row 114 @(112,7999724,"Germany Verbandsliga","Germany Verbandsliga",45375.4375,"Germania OberRoden","RotWeiss Frankfurt",4,0,"H",1.833,4,3.1,1.666,4.333,3.4,-0.75,1.85,1.95,4,1.8,2,0.6659999999999999,-1,-1,0.8500000000000001,-1,0,0)

Set-Row 115 @(113,7999723,"Germany Verbandsliga","Germany Verbandsliga",45375.45833333334,"FC Astoria Walldorf II","FC Muhlhausen 1927",3,1,"H",2.5,3.8,2.2,2.5,3.8,2.2,0,2.025,1.775,3.75,1.85,1.95,1.5,-1,-1,1.025,-1,0.425,-0.5)

Set-Row 116 @(114,7999726,"Germany Verbandsliga","Germany Verbandsliga",45375.47916666666,"SG Andernach","VfB Wissen",3,2,"H",2.5,3.75,2.25,2.625,3.6,2.2,0.25,1.8,2,3.75,1.95,1.85,1.625,-1,-1,0.8,-1,0.95,-1)
